$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 619; everything that was on/after row 619
# shifts down by two (old 619 -> 621, old 620 -> 622, ... old 662 -> 664).
$ws.Rows("619:620").Insert()

# --- New row 619 ---
$ws.Range("A619").Value = 5
$ws.Range("B619").Value = "Macroferia Regional de Talca"
$ws.Range("C619").Value = "Maule"
$ws.Range("D619").Value = 44746
$ws.Range("E619").Value = 7
$ws.Range("F619").Value = 100112020
$ws.Range("G619").Value = "Tomate"
$ws.Range("H619").Value = "Larga vida"
$ws.Range("I619").Value = "Primera"
$ws.Range("J619").Value = 2500
$ws.Range("K619").Value = 10000
$ws.Range("L619").Value = 10000
$ws.Range("M619").Value = 10000
$ws.Range("N619").Value = '$/bandeja 18 kilos'
$ws.Range("O619").Value = "Región de Arica y Parinacota"
$ws.Range("P619").Value = 556
$ws.Range("Q619").Value = 18
$ws.Range("R619").Value = "Hortaliza"

# --- New row 620 ---
$ws.Range("A620").Value = 5
$ws.Range("B620").Value = "Macroferia Regional de Talca"
$ws.Range("C620").Value = "Maule"
$ws.Range("D620").Value = 44746
$ws.Range("E620").Value = 7
$ws.Range("F620").Value = 100112020
$ws.Range("G620").Value = "Tomate"
$ws.Range("H620").Value = "Larga vida"
$ws.Range("I620").Value = "Primera"
$ws.Range("J620").Value = 2000
$ws.Range("K620").Value = 5000
$ws.Range("L620").Value = 5000
$ws.Range("M620").Value = 5000
$ws.Range("N620").Value = '$/caja 10 kilos'
$ws.Range("O620").Value = "Región de Arica y Parinacota"
$ws.Range("P620").Value = 500
$ws.Range("Q620").Value = 10
$ws.Range("R620").Value = "Hortaliza"
